$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Laborations_dagbok - add new diary entries (23/2 -> 1/3-2018) and move the
# "Total tid" summary row down to make room for them.
# ---------------------------------------------------------------------------

# Existing row 10 (23/2-2018) used to read "Player Implementation" - this task
# is now finished, so relabel it and add a brand new row restating that the
# player implementation itself is done.
$ws.Range("B10").Value = "MovingObject done"

# New diary rows 11-15
$ws.Range("A11").Value = "23/2-2018"
$ws.Range("B11").Value = "Player Implementation"
$ws.Range("C11").Value = 2
$ws.Range("D11").Value = 0

$ws.Range("A12").Value = "24/2-2018"
$ws.Range("B12").Value = "sf::Clock Problems"
$ws.Range("C12").Value = 1
$ws.Range("D12").Value = 0

$ws.Range("A13").Value = "1/3-2018"
$ws.Range("B13").Value = "Enemy implementation"
$ws.Range("C13").Value = 0
$ws.Range("D13").Value = 30
# Row 13 previously held the bold "Total tid" label - clear that formatting
# now that the row holds a normal diary entry instead.
$ws.Range("B13").Font.Bold = $false

$ws.Range("A14").Value = "1/3-2018"
$ws.Range("B14").Value = "Player basic completed"
$ws.Range("C14").Value = 0
$ws.Range("D14").Value = 30

$ws.Range("A15").Value = "1/3-2018"
$ws.Range("B15").Value = "Game completed, transformable and overloading"
$ws.Range("C15").Value = 3
$ws.Range("D15").Value = 0

# Move the "Total tid" summary down to row 20, referencing the new data range.
$ws.Range("B20").Value = "Total tid"
$ws.Range("B20").Font.Bold = $true
$ws.Range("C20").Formula = "=SUM(C2:C19)+QUOTIENT(SUM(D2:D19),60)"
$ws.Range("D20").Formula = "=MOD(SUM(D2:D19),60)"

# Column B needs to widen to fit the longer new task descriptions.
$ws.Columns.Item(2).ColumnWidth = 43.67

# Restore the view: scrolled so row 4 is at the top, with B17 selected.
$win = $wb.Windows.Item(1)
$win.ScrollRow = 4
$win.ScrollColumn = 1
$ws.Range("B17").Select() | Out-Null
